# Player.xlsx edit: rename the "View" field to "Cache", and set the
# field's default value to FALSE everywhere it appears.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "View" column header to "Cache" on every sheet that has it ---
# "Property" sheet keeps the field in column F; every "Record_*" sheet keeps
# it in column G (row 1 is the header row).
$ws = $wb.Worksheets.Item("Property")
$ws.Range("F1").Value = "Cache"

$recordSheets = @(
    "Record_PlayerViewItem",
    "Record_PlayerHero",
    "Record_BagEquipList",
    "Record_BagItemList",
    "Record_CommPropertyValue",
    "Record_EctypeList",
    "Record_DropItemList",
    "Record_SkillTable",
    "Record_TaskMonsterList",
    "Record_TaskList",
    "Record_PvpList",
    "Record_ChatGroup",
    "Record_BuildingList",
    "Record_BuildingProduce"
)

foreach ($name in $recordSheets) {
    $rs = $wb.Worksheets.Item($name)
    $rs.Range("G1").Value = "Cache"
}

# --- 2. Set the field's default value to FALSE ---
# "Property" sheet: every data row (2-80) in column F.
$ws = $wb.Worksheets.Item("Property")
$ws.Range("F2:F80").Value = $false

# Every "Record_*" sheet: the single template/default row (row 2) in column G.
foreach ($name in $recordSheets) {
    $rs = $wb.Worksheets.Item($name)
    $rs.Range("G2").Value = $false
}

# --- 3. Restore the workbook's active sheet / window selection ---
# The committed workbook ends up with "Record_BagItemList" (the 5th tab)
# as the active sheet, with cell G1 selected on it, and cell F1 selected
# on the "Property" sheet.
$ws = $wb.Worksheets.Item("Property")
$ws.Range("F1").Select()

foreach ($name in $recordSheets) {
    $rs = $wb.Worksheets.Item($name)
    $rs.Range("G1").Select()
}

$active = $wb.Worksheets.Item("Record_BagItemList")
$active.Activate()
$active.Range("G1").Select()
